# Applies the corrected Scaling/Sampling feature-selection results
# to both the "final_fail" and "final_gifted" worksheets.
$wb = $excel.ActiveWorkbook

function Set-FeatureRow {
    param($ws, $row, $label, $b, $c, $d, $e, $f, $g, $h, $i, $total)
    $ws.Cells.Item($row, 1).Value = $label
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
    $ws.Cells.Item($row, 10).Value = $total
}

# --- final_fail ---
$ws1 = $wb.Worksheets.Item(1)
Set-FeatureRow $ws1, 2, 'Average grade of assignments', $true, $true, $true, $true, $true, $false, $true, $false, 6
Set-FeatureRow $ws1, 3, 'Submissions (% of course total)', $true, $true, $true, $true, $true, $false, $true, $false, 6
Set-FeatureRow $ws1, 4, 'Clicks (% of course total)', $true, $true, $true, $true, $true, $false, $true, $false, 6
Set-FeatureRow $ws1, 5, 'On/off campus click ratio', $true, $true, $true, $true, $true, $false, $true, $false, 6
Set-FeatureRow $ws1, 6, 'Days with no interaction', $true, $true, $false, $true, $true, $false, $true, $false, 5
Set-FeatureRow $ws1, 7, 'Clicks on campus', $false, $true, $true, $true, $true, $false, $true, $false, 5
Set-FeatureRow $ws1, 8, 'Largest period of inactivity (h)', $true, $true, $true, $false, $true, $false, $true, $false, 5
Set-FeatureRow $ws1, 9, 'Total time online (min)', $true, $true, $false, $true, $true, $false, $true, $false, 5
Set-FeatureRow $ws1, 10, 'Start of Session 1 (%)', $true, $true, $false, $true, $true, $false, $true, $false, 5
Set-FeatureRow $ws1, 11, 'Clicks per session', $true, $true, $false, $true, $true, $false, $true, $false, 5
Set-FeatureRow $ws1, 12, 'Number of days', $true, $true, $false, $true, $true, $false, $true, $false, 5
Set-FeatureRow $ws1, 13, 'Resources viewed', $true, $true, $true, $false, $true, $false, $true, $false, 5
Set-FeatureRow $ws1, 14, 'Number of clicks', $true, $true, $false, $false, $true, $false, $true, $false, 4
Set-FeatureRow $ws1, 15, 'Clicks per day', $false, $true, $true, $false, $true, $false, $true, $false, 4
Set-FeatureRow $ws1, 16, 'Clicks on folder', $true, $true, $true, $false, $false, $false, $true, $false, 4
Set-FeatureRow $ws1, 17, 'Average session duration (min)', $true, $true, $false, $false, $true, $false, $true, $false, 4
Set-FeatureRow $ws1, 18, 'Start of Session 3 (%)', $true, $true, $false, $false, $true, $false, $true, $false, 4
Set-FeatureRow $ws1, 19, 'Start of Session 7 (%)', $true, $true, $false, $false, $true, $false, $true, $false, 4
Set-FeatureRow $ws1, 20, 'Links viewed', $true, $true, $false, $false, $false, $false, $true, $false, 3
Set-FeatureRow $ws1, 21, 'Days with no interaction (%)', $false, $true, $false, $false, $true, $false, $true, $false, 3
Set-FeatureRow $ws1, 22, 'Start of Session 2 (%)', $false, $true, $false, $false, $true, $false, $true, $false, 3
Set-FeatureRow $ws1, 23, 'Start of Session 4 (%)', $false, $true, $false, $false, $true, $false, $true, $false, 3
Set-FeatureRow $ws1, 24, 'Assignments viewed', $true, $true, $false, $false, $false, $false, $true, $false, 3
Set-FeatureRow $ws1, 25, 'Files downloaded', $false, $false, $true, $false, $false, $false, $true, $false, 2
Set-FeatureRow $ws1, 26, 'Start of Session 5 (%)', $false, $true, $false, $false, $false, $false, $true, $false, 2
Set-FeatureRow $ws1, 27, 'Start of Session 6 (%)', $false, $true, $false, $false, $false, $false, $true, $false, 2
Set-FeatureRow $ws1, 28, 'Start of Session 10 (%)', $false, $true, $false, $false, $false, $false, $true, $false, 2
Set-FeatureRow $ws1, 29, 'Discussions viewed', $false, $false, $true, $false, $false, $false, $true, $false, 2
Set-FeatureRow $ws1, 30, 'Assignments submitted', $false, $false, $true, $false, $false, $false, $true, $false, 2
Set-FeatureRow $ws1, 31, 'Clicks on course', $false, $true, $false, $false, $false, $false, $true, $false, 2
Set-FeatureRow $ws1, 32, 'Number of sessions', $false, $false, $true, $false, $false, $false, $true, $false, 2
Set-FeatureRow $ws1, 33, 'Clicks on forum', $false, $true, $false, $false, $false, $false, $true, $false, 2
Set-FeatureRow $ws1, 34, 'Forum posts', $false, $false, $true, $false, $false, $false, $true, $false, 2
Set-FeatureRow $ws1, 35, 'Start of Session 9 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1
Set-FeatureRow $ws1, 36, 'Quizzes started', $false, $false, $false, $false, $false, $false, $true, $false, 1
Set-FeatureRow $ws1, 37, 'Start of Session 8 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1

# --- final_gifted ---
$ws2 = $wb.Worksheets.Item(2)
Set-FeatureRow $ws2, 2, 'Clicks (% of course total)', $true, $true, $true, $true, $true, $false, $true, $false, 6
Set-FeatureRow $ws2, 3, 'Average grade of assignments', $true, $true, $true, $true, $true, $false, $true, $false, 6
Set-FeatureRow $ws2, 4, 'Total time online (min)', $true, $true, $true, $true, $true, $false, $true, $false, 6
Set-FeatureRow $ws2, 5, 'Average session duration (min)', $true, $true, $true, $true, $true, $false, $true, $false, 6
Set-FeatureRow $ws2, 6, 'Start of Session 1 (%)', $true, $true, $true, $true, $true, $false, $true, $false, 6
Set-FeatureRow $ws2, 7, 'Largest period of inactivity (h)', $true, $true, $true, $false, $true, $false, $true, $false, 5
Set-FeatureRow $ws2, 8, 'Days with no interaction', $true, $true, $false, $true, $true, $false, $true, $false, 5
Set-FeatureRow $ws2, 9, 'On/off campus click ratio', $true, $true, $false, $true, $true, $false, $true, $false, 5
Set-FeatureRow $ws2, 10, 'Clicks per session', $true, $true, $true, $false, $true, $false, $true, $false, 5
Set-FeatureRow $ws2, 11, 'Assignments viewed', $true, $true, $false, $false, $true, $false, $true, $false, 4
Set-FeatureRow $ws2, 12, 'Resources viewed', $true, $false, $false, $true, $true, $false, $true, $false, 4
Set-FeatureRow $ws2, 13, 'Number of days', $true, $true, $false, $false, $true, $false, $true, $false, 4
Set-FeatureRow $ws2, 14, 'Clicks per day', $true, $false, $false, $false, $true, $false, $true, $false, 3
Set-FeatureRow $ws2, 15, 'Assignments submitted', $true, $false, $true, $false, $false, $false, $true, $false, 3
Set-FeatureRow $ws2, 16, 'Start of Session 6 (%)', $true, $false, $false, $false, $true, $false, $true, $false, 3
Set-FeatureRow $ws2, 17, 'Days with no interaction (%)', $true, $false, $false, $false, $true, $false, $true, $false, 3
Set-FeatureRow $ws2, 18, 'Start of Session 2 (%)', $true, $false, $false, $false, $true, $false, $true, $false, 3
Set-FeatureRow $ws2, 19, 'Submissions (% of course total)', $false, $false, $false, $false, $true, $false, $true, $false, 2
Set-FeatureRow $ws2, 20, 'Files downloaded', $false, $false, $true, $false, $false, $false, $true, $false, 2
Set-FeatureRow $ws2, 21, 'Quizzes started', $false, $false, $true, $false, $false, $false, $true, $false, 2
Set-FeatureRow $ws2, 22, 'Number of clicks', $false, $false, $false, $false, $true, $false, $true, $false, 2
Set-FeatureRow $ws2, 23, 'Clicks on folder', $false, $false, $true, $false, $false, $false, $true, $false, 2
Set-FeatureRow $ws2, 24, 'Clicks on course', $true, $false, $false, $false, $false, $false, $true, $false, 2
Set-FeatureRow $ws2, 25, 'Start of Session 7 (%)', $false, $false, $false, $false, $true, $false, $true, $false, 2
Set-FeatureRow $ws2, 26, 'Start of Session 5 (%)', $false, $false, $false, $false, $true, $false, $true, $false, 2
Set-FeatureRow $ws2, 27, 'Start of Session 4 (%)', $false, $false, $false, $false, $true, $false, $true, $false, 2
Set-FeatureRow $ws2, 28, 'Start of Session 3 (%)', $false, $false, $false, $false, $true, $false, $true, $false, 2
Set-FeatureRow $ws2, 29, 'Clicks on campus', $false, $false, $false, $false, $true, $false, $true, $false, 2
Set-FeatureRow $ws2, 30, 'Links viewed', $false, $false, $true, $false, $false, $false, $true, $false, 2
Set-FeatureRow $ws2, 31, 'Discussions viewed', $false, $false, $false, $false, $false, $false, $true, $false, 1
Set-FeatureRow $ws2, 32, 'Number of sessions', $false, $false, $false, $false, $false, $false, $true, $false, 1
Set-FeatureRow $ws2, 33, 'Forum posts', $false, $false, $false, $false, $false, $false, $true, $false, 1
Set-FeatureRow $ws2, 34, 'Clicks on forum', $false, $false, $false, $false, $false, $false, $true, $false, 1
Set-FeatureRow $ws2, 35, 'Start of Session 10 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1
Set-FeatureRow $ws2, 36, 'Start of Session 9 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1
Set-FeatureRow $ws2, 37, 'Start of Session 8 (%)', $false, $false, $false, $false, $false, $false, $true, $false, 1

Write-Output "Feature selection rows updated."
